# US 3.3 commit files
# Apply content/structure changes to "About" (sheet1) and "DR" (sheet2) sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("DR")

# ---------------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------------

# Year of the source table moved from 2016 to 2015
$ws1.Range("B4").Value = 2015

# Source link: was the EPA 2016 SCC TSD, now the 2015 OMB/whitehouse SCC TSD,
# referenced via a hyperlink to the interagency working-group circular.
$ws1.Hyperlinks.Add($ws1.Range("B6"), "http://www.whitehouse.gov/omb/circulars_a094#8", "8", "", "http://www.whitehouse.gov/omb/circulars_a094#8") | Out-Null
$ws1.Range("B6").Value = "https://www.whitehouse.gov/sites/default/files/omb/inforeg/scc-tsd-final-july-2015.pdf"

# Notes section: rows 10-12 now describe the "Annual Perc" input instead of
# the old SCC discount-rate rationale text, which is pushed down (with a
# blank separator row 15) into rows 16-18.
$ws1.Range("A10").Value = "This is the annual percentage rate by which future savings (e.g. fuel cost savings) are discounted when"
$ws1.Range("A11").Value = "making price-driven purchasing decisions in the current year. The value used should be one that is"
$ws1.Range("A12").Value = "reasonable for people who are looking to buy fuel-consuming capital equipment, such as industrial"
$ws1.Range("A13").Value = "equipment or building components. The model works in real dollars, so this rate should be the growth"
$ws1.Range("A14").Value = "in real value, not the growth in nominal value plus real value."

$ws1.Range("A16").Value = "We choose to use a 3% discount rate here, for consistency with the 3% rate used for the central estimate"
$ws1.Range("A17").Value = "of Social Cost of Carbon (in the SCoC variable), as well as the discount rate built into the health"
$ws1.Range("A18").Value = "damages values in the SCoHIbP Social Cost of Health Impacts by Pollutant variable."

# ---------------------------------------------------------------------------
# Sheet "DR"
# ---------------------------------------------------------------------------

$ws2.Range("B1").Value = "Annual Perc (dimensionless)"
$ws2.Range("B1").WrapText = $true
$ws2.Rows.Item(1).RowHeight = 30

# ---------------------------------------------------------------------------
# Selections (match the saved view state from the authored workbook)
# ---------------------------------------------------------------------------

$excel.Goto($ws2.Range("B1"))
$excel.Goto($ws1.Range("A16:A18"))
